# Edit script: reshape "Multiplying Fractions" deck per commit
# "[Presentation] Presentation created through GPT prompt"

$p = $ppt.ActivePresentation
$layouts = $p.SlideMaster.CustomLayouts

# Handy layout handles (match slideLayoutN.xml order / names)
$lytTitleSlide   = $layouts.Item(1)   # Title Slide
$lytTitleContent = $layouts.Item(2)   # Title and Content
$lytSectionHdr   = $layouts.Item(3)   # Section Header
$lytTwoContent   = $layouts.Item(4)   # Two Content
$lytComparison   = $layouts.Item(5)   # Comparison
$lytTitleOnly    = $layouts.Item(6)   # Title Only
$lytBlank        = $layouts.Item(7)   # Blank
$lytContentCap   = $layouts.Item(8)   # Content with Caption
$lytPictureCap   = $layouts.Item(9)   # Picture with Caption

# ---------------------------------------------------------------------------
# Slide 1 ("Multiplying Fractions" title slide): clear the subtitle text.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = ""

# ---------------------------------------------------------------------------
# Slide 2 ("Introduction" -> "What is Fraction Multiplication?")
# Switch to the Section Header layout.
# ---------------------------------------------------------------------------
$p.Slides.Item(2).Delete()
$s2 = $p.Slides.AddSlide(2, $lytSectionHdr)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "What is Fraction Multiplication?"
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "A Mathematical Operation Applied to Fractions"

# ---------------------------------------------------------------------------
# Slide 3 ("What are Fractions?" -> "Key Concept") : layout unchanged.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Key Concept"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "When we multiply fractions, we are finding a part of a part."

# ---------------------------------------------------------------------------
# Slide 4 ("Basics of Multiplying Fractions" -> "Multiplication vs. Addition")
# Switch to the Comparison layout (5 placeholders).
# ---------------------------------------------------------------------------
$p.Slides.Item(4).Delete()
$s4 = $p.Slides.AddSlide(4, $lytComparison)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Multiplication vs. Addition"
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "Multiplying Fractions"
$s4.Shapes.Item(3).TextFrame.TextRange.Text = "Multiply the numerators and the denominators."
$s4.Shapes.Item(4).TextFrame.TextRange.Text = "Adding Fractions"
$s4.Shapes.Item(5).TextFrame.TextRange.Text = "Find a common denominator before adding numerators."

# ---------------------------------------------------------------------------
# Slide 5 ("Example 1" -> "Steps to Multiply Fractions")
# Switch to the Two Content layout.
# ---------------------------------------------------------------------------
$p.Slides.Item(5).Delete()
$s5 = $p.Slides.AddSlide(5, $lytTwoContent)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Steps to Multiply Fractions"
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "1. Multiply the numerators."
$s5.Shapes.Item(3).TextFrame.TextRange.Text = "2. Multiply the denominators."

# ---------------------------------------------------------------------------
# Slide 6 ("Example 2" -> "Example Problem")
# Switch to the Content with Caption layout.
# ---------------------------------------------------------------------------
$p.Slides.Item(6).Delete()
$s6 = $p.Slides.AddSlide(6, $lytContentCap)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Example Problem"
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Let's multiply 3/4 by 2/3."
$s6.Shapes.Item(3).TextFrame.TextRange.Text = "Simplifying the Product"

# ---------------------------------------------------------------------------
# Slide 7 ("Practice Problem" -> "Visual Representation")
# Switch to the Picture with Caption layout; leave picture + caption empty.
# ---------------------------------------------------------------------------
$p.Slides.Item(7).Delete()
$s7 = $p.Slides.AddSlide(7, $lytPictureCap)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Visual Representation"

# ---------------------------------------------------------------------------
# Slide 8 ("Review and Q&A" -> "Practice Time!")
# Switch to the Title Only layout (drop the body placeholder).
# ---------------------------------------------------------------------------
$p.Slides.Item(8).Delete()
$s8 = $p.Slides.AddSlide(8, $lytTitleOnly)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Practice Time!"

# ---------------------------------------------------------------------------
# Slide 9 ("Conclusion" -> completely empty slide)
# Switch to the Blank layout (no placeholders at all).
# ---------------------------------------------------------------------------
$p.Slides.Item(9).Delete()
$s9 = $p.Slides.AddSlide(9, $lytBlank)

# ---------------------------------------------------------------------------
# Slide 10 (new "Summary" slide) appended at the end.
# ---------------------------------------------------------------------------
$s10 = $p.Slides.AddSlide(10, $lytTitleContent)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Summary"
$s10.Shapes.Item(2).TextFrame.TextRange.Text = "We've learned the steps of multiplying fractions and seen an example. Now it's time to practice on your own!"
